$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / link / percent-style cells: safe to assign directly ---
$ws.Range("D2").Value = '68.738.39'
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").Value = '3.735.08'
$ws.Range("E3").Value = '  -1.09%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("E6").Value = '  -4.37%  '
$ws.Range("D7").Value = '3.732.10'
$ws.Range("E7").Value = '  -1.18%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +1.53%  '
$ws.Range("E10").Value = '  +3.21%  '
$ws.Range("E11").Value = '  +3.11%  '
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").Value = '4.361.42'
$ws.Range("E15").Value = '  -1.14%  '
$ws.Range("D16").Value = '3.737.88'
$ws.Range("E16").Value = '  -1.50%  '
$ws.Range("D17").Value = '68.730.10'
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("E20").Value = '  +4.45%  '
$ws.Range("E21").Value = '  +1.25%  '
$ws.Range("E22").Value = '  +10.85%  '
$ws.Range("E23").Value = '  -2.52%  '
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  -4.71%  '
$ws.Range("E26").Value = '  -2.43%  '
$ws.Range("E27").Value = '  +1.11%  '
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E32").Value = '  +2.76%  '
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").Value = '3.881.48'
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("D35").Value = '3.668.49'
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("E38").Value = '  +0.63%  '
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("E40").Value = '  -1.41%  '
$ws.Range("E41").Value = '  -1.32%  '
$ws.Range("E42").Value = '  -3.22%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("E45").Value = '  -1.42%  '
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("E49").Value = '  +1.47%  '
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("D51").Value = '2.742.06'
$ws.Range("E51").Value = '  -3.31%  '

# --- Price cells whose new text would otherwise be auto-parsed as a number: ---
# force Text format, assign, then restore the default "Normal" style so no
# extra formatting is left behind (matches original plain inlineStr cells).
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '601.80'
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '167.56'
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.535'
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '6.38'
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '17.25'
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '495.23'
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '10.03'
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.722'
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '84.88'
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.30'
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '12.39'
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '10.08'
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '1.01'
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '5.81'
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.133'
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.324'
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '433.46'
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '49.19'
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '2.88'
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.97'
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '8.41'
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '40.68'
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '140.98'
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.0352'
$cell.Style = "Normal"
